$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap LEG1_DIRECTION (AL2) and LEG2_DIRECTION (AM2) values: "P" <-> "R"
$ws.Range("AL2").Value = "R"
$ws.Range("AM2").Value = "P"

# Scroll the view so column AB is the left-most visible column and select AM2,
# matching the author's final cursor/scroll position when the file was saved.
$excel.ActiveWindow.ScrollColumn = 28
$ws.Range("AM2").Select()
